# The document's second table (the "signalement levels" nomenclature
# table) was generated from a pandas DataFrame with a lot of spurious
# "Unnamed: N" columns and blank trailing rows. This edit trims it back
# down to the 5 meaningful columns (Code / Libellé niveau 1 / Libellé
# niveau 2 / Description / Commentaire) and the 3 meaningful rows
# (header + STANDARD + ATTENTION), then widens the remaining columns
# from 665 dxa to 1728 dxa so they fill the same overall table width.

$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# Drop the empty trailing rows (rows 4..15), working from the bottom up
# so earlier indices stay valid as rows are removed.
for ($i = $t.Rows.Count; $i -ge 4; $i--) {
    $t.Rows.Item($i).Delete()
}

# Drop the "Unnamed: 5" .. "Unnamed: 12" columns (columns 6..13),
# again working from the right so earlier indices stay valid.
for ($i = $t.Columns.Count; $i -ge 6; $i--) {
    $t.Columns.Item($i).Delete()
}

# Re-widen the 5 remaining columns: 1728 dxa == 86.4 points (20 dxa/pt).
for ($i = 1; $i -le $t.Columns.Count; $i++) {
    $t.Columns.Item($i).Width = 86.4
}
